$d = $word.ActiveDocument

# Insert "Faith" in the middle of the First Name underscores:
# "First Name: ________" -> "First Name: __Faith______"
$rFirst = $d.Content
$found = $rFirst.Find.Execute("First Name: __", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rFirst.Collapse(0)
    $rFirst.InsertAfter("Faith")
}

# Insert "Johnson" in the middle of the Last Name underscores:
# "Last Name: ________" -> "Last Name: __Johnson______"
$rLast = $d.Content
$found = $rLast.Find.Execute("Last Name: __", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rLast.Collapse(0)
    $rLast.InsertAfter("Johnson")
}
